$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.334.79"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.687.24"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "680.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.441"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "4.312.84"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "3.695.59"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "69.317.09"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "3.837.09"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "3.677.86"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "171.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.942"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000274"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.88%  "
